# Update the route-section reference codes in column A (Sheet1, rows 2-32).
# Old scheme:  TR/8350/12AB/00/<1|3>/<date>
# New scheme:  TR/8350/12AB/2021/<10.01|20.01>/<date>
# (see commit message: "changed versioning of route sections in all examples")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRefs = @(
    "TR/8350/12AB/2021/10.01/2021-12-01",
    "TR/8350/12AB/2021/10.01/2021-12-02",
    "TR/8350/12AB/2021/20.01/2021-12-03",
    "TR/8350/12AB/2021/20.01/2021-12-04",
    "TR/8350/12AB/2021/20.01/2021-12-05",
    "TR/8350/12AB/2021/10.01/2021-12-06",
    "TR/8350/12AB/2021/10.01/2021-12-07",
    "TR/8350/12AB/2021/10.01/2021-12-08",
    "TR/8350/12AB/2021/10.01/2021-12-09",
    "TR/8350/12AB/2021/20.01/2021-12-10",
    "TR/8350/12AB/2021/20.01/2021-12-11",
    "TR/8350/12AB/2021/20.01/2021-12-12",
    "TR/8350/12AB/2021/10.01/2021-12-13",
    "TR/8350/12AB/2021/10.01/2021-12-14",
    "TR/8350/12AB/2021/10.01/2021-12-15",
    "TR/8350/12AB/2021/10.01/2021-12-16",
    "TR/8350/12AB/2021/20.01/2021-12-17",
    "TR/8350/12AB/2021/20.01/2021-12-18",
    "TR/8350/12AB/2021/20.01/2021-12-19",
    "TR/8350/12AB/2021/10.01/2021-12-20",
    "TR/8350/12AB/2021/10.01/2021-12-21",
    "TR/8350/12AB/2021/10.01/2021-12-22",
    "TR/8350/12AB/2021/10.01/2021-12-23",
    "TR/8350/12AB/2021/20.01/2021-12-24",
    "TR/8350/12AB/2021/20.01/2021-12-25",
    "TR/8350/12AB/2021/20.01/2021-12-26",
    "TR/8350/12AB/2021/10.01/2021-12-27",
    "TR/8350/12AB/2021/10.01/2021-12-28",
    "TR/8350/12AB/2021/10.01/2021-12-29",
    "TR/8350/12AB/2021/10.01/2021-12-30",
    "TR/8350/12AB/2021/20.01/2021-12-31"
)

for ($i = 0; $i -lt $newRefs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newRefs[$i]
}

# Match the author's final UI view/selection state recorded in the workbook:
# zoom reset back to 100% (no custom zoomScale) and a new active cell.
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("G5").Select() | Out-Null
